# Updated readme and backlog
#
# The deck's "automatically updated" date placeholder (datetimeFigureOut
# field) on the slide master and every slide layout rolls from
# 29/03/2023 to 30/03/2023, and the slow auto-advance transition timing
# moves from 5 seconds (5000 ms) to 30 seconds (30000 ms) everywhere it
# appears: the slide master, every slide layout, and all three slides.

$p = $ppt.ActivePresentation

# ppPlaceholderDate
$ppPlaceholderDate = 16

function Update-DatePlaceholderText {
    param($shapes, $newText)

    for ($i = 1; $i -le $shapes.Count; $i++) {
        $shp = $shapes.Item($i)
        $isDatePlaceholder = $false
        try {
            if ($shp.PlaceholderFormat.Type -eq $ppPlaceholderDate) {
                $isDatePlaceholder = $true
            }
        } catch {
            $isDatePlaceholder = $false
        }

        if ($isDatePlaceholder) {
            try {
                $shp.TextFrame.TextRange.Text = $newText
            } catch {
            }
        }
    }
}

function Update-LayoutTransitionTiming {
    param($transitionHolder, $seconds)

    try {
        $transitionHolder.SlideShowTransition.AdvanceTime = $seconds
    } catch {
    }
    try {
        $transitionHolder.SlideShowTransition.AdvanceOnTime = -1
    } catch {
    }
}

$newDateText = "30/03/2023"
$newAdvanceSeconds = 30

# --- Slide master: date placeholder + transition timing ---
$master = $p.SlideMaster
Update-DatePlaceholderText $master.Shapes $newDateText
Update-LayoutTransitionTiming $master $newAdvanceSeconds

# --- Every slide layout: date placeholder + transition timing ---
for ($L = 1; $L -le $master.CustomLayouts.Count; $L++) {
    $layout = $master.CustomLayouts.Item($L)
    Update-DatePlaceholderText $layout.Shapes $newDateText
    Update-LayoutTransitionTiming $layout $newAdvanceSeconds
}

# --- Every slide: transition timing (auto-advance after 30s) ---
for ($S = 1; $S -le $p.Slides.Count; $S++) {
    $slide = $p.Slides.Item($S)
    $slide.SlideShowTransition.AdvanceTime = $newAdvanceSeconds
}
